# "regenerate orders with updates distance/sizes"
#
# The trial-order sheet encodes experiment parameters (viewing Distance and
# stimulus Size) directly inside text tokens such as "Face02_D64_S30",
# "Fixation_D64_l.png", "D64", "S30", etc. This regenerates the order with
# new distance/size values:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# (S25 / S20 are untouched.)
#
# This is a pure text substitution applied to every string cell in the
# sheet - the row/column layout and the cell->content mapping are
# otherwise unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            $newVal = $val.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
